$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Shared-string text update -------------------------------------------
# Q1 (merged Q1:U1, 3rd sprint header) keeps pointing at the same shared
# string slot, but its wording changes ("Vorarbeiten zu Spieler entfernen"
# replaces the old goal text).
$ws.Range("Q1").Value = "Sprintziel: Spieler anlegen, Vorarbeiten zu Spieler entfernen"

# --- 5th sprint column (W:AB) "Occupations" ------------------------------
# New occupation/topic entries added next to each daily-scrum date (W col)
# for the 5th sprint, plus the updated sprint-goal header in X1.
$ws.Range("X5").Value = "Planung zur weiteren Vorgehensweise"
$ws.Range("X1").Value = "Sprintziel: Spieler entfernen, Spiel anlegen, Positionsverwaltung, Login"
$ws.Range("X7").Value = "Besprechung zu Spieler entfernen"
$ws.Range("X9").Value = "Besprechung zur Positionsverwaltung"
$ws.Range("X11").Value = "Besprechung zu Login und Spiel anlegen"
$ws.Range("X15").Value = "Planung zur Präsentation für 23.5."
$ws.Range("X19").Value = "Sprint Review"
$ws.Range("X3").Value = "Sprintplanung"
$ws.Range("X13").Value = "Besprechung Akzeptanzkriterien"
$ws.Range("X17").Value = "Abschluss Akzeptanzkriterien"

# Widen column AB so the new occupation text fits.
$ws.Range("AB1").ColumnWidth = 17.66

# --- View state -----------------------------------------------------------
# Scroll the visible pane so column S is at the left edge, then move the
# selection to AG13.
$excel.Goto($ws.Range("S1"), $true)
$ws.Range("AG13").Select()
